$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $value
    $c.Style = "Normal"
}

Set-TextValue 'D2' '42.556.87'
$ws.Range('E2').Value = '  -0.98%  '
Set-TextValue 'D3' '2.532.89'
$ws.Range('E3').Value = '  -1.67%  '
$ws.Range('E4').Value = '  -0.05%  '
Set-TextValue 'D5' '305.14'
$ws.Range('E5').Value = '  +0.72%  '
Set-TextValue 'D6' '97.42'
$ws.Range('E6').Value = '  -0.11%  '
Set-TextValue 'D7' '0.592'
$ws.Range('E7').Value = '  +2.59%  '
$ws.Range('E8').Value = '  +0.08%  '
$ws.Range('E9').Value = '  -2.18%  '
Set-TextValue 'D10' '36.81'
$ws.Range('E10').Value = '  +0.88%  '
Set-TextValue 'D11' '0.0812'
$ws.Range('E11').Value = '  +0.36%  '
$ws.Range('E12').Value = '  +0.28%  '
Set-TextValue 'D13' '0.113'
$ws.Range('E13').Value = '  -1.35%  '
Set-TextValue 'D14' '2.920.41'
$ws.Range('E14').Value = '  -1.56%  '
Set-TextValue 'D15' '2.520.54'
$ws.Range('E15').Value = '  -3.17%  '
Set-TextValue 'D16' '15.26'
$ws.Range('E16').Value = '  +5.94%  '
Set-TextValue 'D17' '0.865'
$ws.Range('E17').Value = '  -2.59%  '
Set-TextValue 'D18' '42.584.12'
$ws.Range('E18').Value = '  -1.06%  '
Set-TextValue 'D19' '12.99'
$ws.Range('E19').Value = '  +0.21%  '
Set-TextValue 'D20' '0.0₃0977'
$ws.Range('E20').Value = '  -1.95%  '
Set-TextValue 'D21' '6.48'
$ws.Range('E21').Value = '  -2.59%  '
Set-TextValue 'D22' '71.18'
$ws.Range('E22').Value = '  -1.15%  '
Set-TextValue 'D23' '251.78'
$ws.Range('E23').Value = '  -1.30%  '
Set-TextValue 'D24' '2.93'
$ws.Range('E24').Value = '  -1.38%  '
Set-TextValue 'D25' '2.03'
$ws.Range('E25').Value = '  -4.66%  '
Set-TextValue 'D26' '27.05'
$ws.Range('E26').Value = '  -5.90%  '
$ws.Range('E27').Value = '  +0.10%  '
Set-TextValue 'D28' '2.34'
$ws.Range('E28').Value = '  +10.03%  '
Set-TextValue 'D29' '10.34'
$ws.Range('E29').Value = '  +0.68%  '
$ws.Range('E30').Value = '  +1.16%  '
Set-TextValue 'D31' '5.99'
$ws.Range('E31').Value = '  -1.32%  '
Set-TextValue 'D32' '155.50'
$ws.Range('E32').Value = '  -0.05%  '
Set-TextValue 'D33' '3.32'
$ws.Range('E33').Value = '  -2.32%  '
Set-TextValue 'D34' '0.0792'
$ws.Range('E34').Value = '  -2.47%  '
$ws.Range('B35').Value = 'ARBITRUM'
$ws.Range('C35').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue 'D35' '2.08'
$ws.Range('E35').Value = '  -5.13%  '
$ws.Range('B36').Value = 'Celestia'
$ws.Range('C36').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
Set-TextValue 'D36' '18.70'
$ws.Range('E36').Value = '  +1.68%  '
$ws.Range('E37').Value = '  -4.66%  '
Set-TextValue 'D38' '0.116'
$ws.Range('E38').Value = '  +1.56%  '
$ws.Range('E39').Value = '  +0.09%  '
Set-TextValue 'D40' '24.13'
$ws.Range('E40').Value = '  +1.44%  '
$ws.Range('E41').Value = '  -1.04%  '
Set-TextValue 'D42' '3.87'
$ws.Range('E42').Value = '  -0.47%  '
$ws.Range('B43').Value = 'FirstDigitalUSD'
$ws.Range('C43').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextValue 'D43' '0.998'
$ws.Range('E43').Value = '  +0.08%  '
$ws.Range('B44').Value = 'VeChain'
$ws.Range('C44').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue 'D44' '0.0301'
$ws.Range('E44').Value = '  -3.27%  '
$ws.Range('B45').Value = 'ApeXProtocol'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
Set-TextValue 'D45' '2.00'
$ws.Range('E45').Value = '  -1.81%  '
Set-TextValue 'D46' '2.043.70'
$ws.Range('E46').Value = '  -1.28%  '
Set-TextValue 'D47' '85.00'
$ws.Range('E47').Value = '  -0.65%  '
Set-TextValue 'D48' '8.98'
$ws.Range('E48').Value = '  -3.33%  '
Set-TextValue 'D49' '2.778.11'
$ws.Range('E49').Value = '  -1.56%  '
$ws.Range('B50').Value = 'Algorand'
$ws.Range('C50').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue 'D50' '0.190'
$ws.Range('E50').Value = '  -0.91%  '
$ws.Range('B51').Value = 'Aave'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue 'D51' '102.39'
$ws.Range('E51').Value = '  -4.04%  '
